# Cleaning script for raw adcap_score data from the experts:
# convert the numeric 0/0.5/1-ish scores on a few sheets into the
# proper categorical labels (none/low/medium/high), fix a couple of
# mis-scored cells, and leave the selection/active-tab where the
# editor last left off.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "biomass removal" -- every adcap_score in D2:D15 was
# raw/unscored data; recode them all to "none".
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("biomass removal")
$ws1.Activate()

$ws1.Range("D2").Value  = "none"
$ws1.Range("D3").Value  = "none"
$ws1.Range("D4").Value  = "none"
$ws1.Range("D5").Value  = "none"
$ws1.Range("D6").Value  = "none"
$ws1.Range("D7").Value  = "none"
$ws1.Range("D8").Value  = "none"
$ws1.Range("D9").Value  = "none"
$ws1.Range("D10").Value = "none"
$ws1.Range("D11").Value = "none"
$ws1.Range("D12").Value = "none"
$ws1.Range("D13").Value = "none"
$ws1.Range("D14").Value = "none"
$ws1.Range("D15").Value = "none"

$ws1.Range("D2:D8").Select()

# ---------------------------------------------------------------
# Sheet 3: "habitat loss + degradation" -- recode the adult
# mobility block (D2:D8) to labels, and clear out the bad
# PLD-exposure scores in D10:D15 back to zero.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("habitat loss + degradation")
$ws3.Activate()

$ws3.Range("D2").Value = "none"
$ws3.Range("D3").Value = "none"
$ws3.Range("D4").Value = "low"
$ws3.Range("D5").Value = "medium"
$ws3.Range("D6").Value = "high"
$ws3.Range("D7").Value = "high"
$ws3.Range("D8").Value = "high"

$ws3.Range("D10").Value = 0
$ws3.Range("D11").Value = 0
$ws3.Range("D12").Value = 0
$ws3.Range("D13").Value = 0
$ws3.Range("D14").Value = 0
$ws3.Range("D15").Value = 0

$ws3.Range("D9:D15").Select()

# ---------------------------------------------------------------
# Sheet 5: "poisons + toxins" -- recode the zone block's D2 and
# the PLD-exposure block D9:D15 to labels, plus a note on D9/E9
# explaining the assumption behind the fix.
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("poisons + toxins")
$ws5.Activate()

$ws5.Range("D2").Value = "none"

$ws5.Range("D9").Value = "high"
$ws5.Range("E9").Value = "assume short PLD means less exposure to poisons and toxins"

$ws5.Range("D10").Value = "medium"
$ws5.Range("D11").Value = "medium"
$ws5.Range("D12").Value = "low"
$ws5.Range("D13").Value = "low"
$ws5.Range("D14").Value = "none"
$ws5.Range("D15").Value = "high"

$ws5.Range("D3").Select()
